# VP_FILTER_SKU.xlsx - refresh the SKU list on "vp_sku_list"
# Replaces the old 5-SKU sample (rows 2-6, then blank styled rows 7-34)
# with the new 77-row SKU list (rows 2-78), matching the latest export.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$skus = @(
    10022921,10046772,10177346,10173048,10022921,10046772,10030631,10134810,
    10177346,10046772,10020378,10020378,10217540,10177346,10127368,10036227,
    10114434,10134810,10046772,10113223,10020378,10022921,10177346,10005480,
    10127368,10134810,10187607,10059535,10134810,10164431,10114434,10219764,
    10134810,10036227,10113223,10127368,10177346,10187607,10114434,10113223,
    10024943,10114434,10219764,10059024,10127368,10164340,10113223,10114434,
    10004408,10001043,10219764,10127368,10255093,10036227,10024943,10114434,
    10114434,10113223,10164340,10059024,10164431,10127368,10024943,10114434,
    10113223,10016287,10232805,10024943,10255093,10059139,10059024,10113223,
    10164340,10255093,10114434,10024943,10004408
)

$firstRow = 2
$lastRow = $firstRow + $skus.Count - 1

for ($i = 0; $i -lt $skus.Count; $i++) {
    $ws.Cells.Item($firstRow + $i, 1).Value = $skus[$i]
}

# Select the refreshed data range, matching the saved workbook state.
$ws.Range("A2:A" + $lastRow).Select() | Out-Null
